# Practica#2 - Cambio de contraseña antes de ingresar al sistema.
# Edits to betali_v2_asistencias/extras/temas.xlsx ("temas" theme table):
#   - Row 10 ("CHyP" theme): color_letra typo'd from "#fff" to "$fff"
#   - Row 11 ("Cool" theme) removed entirely
#   - Selection left on column G (fecha_registro) after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("temas")
$ws.Activate()

# Fix/typo the color_letra value for the "CHyP" theme row (row 10).
$ws.Range("C10").Value = '$fff'

# Remove the "Cool" theme row (row 11) entirely.
$ws.Rows(11).Delete()

# Leave the selection on column G, matching the saved state of the sheet.
[void]$ws.Columns("G").Select()
